$d = $word.ActiveDocument

# Find the paragraph containing the target sentence and remove the whole
# paragraph (including its paragraph mark) so the list collapses cleanly.
$range = $d.Content
$found = $range.Find.Execute("a function that calculates projected income for hosts.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $range.Paragraphs(1)
    $paraRange = $para.Range
    $paraRange.Delete()
}
